# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = -0.6886020151133502
$ws.Range("H2").Value = -1.031801007556675
$ws.Range("I2").Value = -0.8167506297229219
$ws.Range("J2").Value = -0.8167506297229219
$ws.Range("K2").Value = -1.961
$ws.Range("L2").Value = -0.6174433249370278
$ws.Range("U2").Value = 3.08
$ws.Range("V2").Value = 0.1632220455749867
$ws.Range("W2").Value = -1.172989686568367
$ws.Range("X2").Value = 0.05327491537410382
$ws.Range("Y2").Value = -1.226264601942471
$ws.Range("Z2").Value = -1.962917181705808
$ws.Range("AA2").Value = -0.1229312279804526
$ws.Range("AB2").Value = 0.05327012544756154
$ws.Range("AC2").Value = -0.1762013534280142
$ws.Range("AD2").Value = 0.005
$ws.Range("AF2").Value = 0.005
$ws.Range("AG2").Value = -3.075
$ws.Range("AH2").Value = 0.0002649006622516556
$ws.Range("AI2").Value = 0.0008554319931565441
$ws.Range("AJ2").Value = -0.194681861348528
$ws.Range("AK2").Value = -1.112115732368897
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = -0.251
$ws.Range("AN2").Value = -0.001955416503715291
$ws.Range("AO2").ClearContents()
$ws.Range("AP2").Value = 1.202581149784904
$ws.Range("AQ2").Value = 10.33466135458167
$ws.Range("G3").Value = -270.3333333333333
$ws.Range("H3").Value = -410
$ws.Range("I3").Value = -291.6666666666667
$ws.Range("J3").Value = -291.6666666666667
$ws.Range("K3").Value = -1.28
$ws.Range("L3").Value = -213.3333333333333
$ws.Range("U3").Value = 1.3
$ws.Range("V3").Value = 0.8280254777070064
$ws.Range("W3").Value = -2.165820642978003
$ws.Range("X3").Value = 0.05326995922263075
$ws.Range("Y3").Value = -2.219090602200634
$ws.Range("Z3").Value = -0.002386634844868735
$ws.Range("AA3").Value = 0.696101829753381
$ws.Range("AB3").Value = 0.05326995922263075
$ws.Range("AC3").Value = 0.6428318705307502
$ws.Range("AD3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -1.3
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -4.814814814814815
$ws.Range("AK3").Value = -0.8783783783783785
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = -0.205
$ws.Range("AN3").Value = -0
$ws.Range("AO3").ClearContents()
$ws.Range("AP3").Value = 0.7471264367816092
$ws.Range("AQ3").Value = 8.536585365853659
$ws.Range("G4").Value = -0.1782334384858044
$ws.Range("H4").Value = -0.2577287066246057
$ws.Range("I4").Value = -0.2662460567823344
$ws.Range("J4").Value = -0.2662460567823344
$ws.Range("K4").Value = -0.681
$ws.Range("L4").Value = -0.214826498422713
$ws.Range("U4").Value = 1.78
$ws.Range("V4").Value = 0.1028901734104046
$ws.Range("W4").Value = -0.1801587301587302
$ws.Range("X4").Value = 0.05327987152557688
$ws.Range("Y4").Value = -0.2334386016843071
$ws.Range("Z4").Value = 3.537946428571431
$ws.Range("AA4").Value = -0.9419642857142863
$ws.Range("AB4").Value = 0.05327029167249234
$ws.Range("AC4").Value = -0.9952345773867786
$ws.Range("AD4").Value = 0.005
$ws.Range("AF4").Value = 0.005
$ws.Range("AG4").Value = -1.775
$ws.Range("AH4").Value = 0.0002889338341519792
$ws.Range("AI4").Value = 0.001631321370309951
$ws.Range("AJ4").Value = -0.1143317230273752
$ws.Range("AK4").Value = -1.381322957198444
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = -0.046
$ws.Range("AN4").Value = -0.006119951040391677
$ws.Range("AO4").ClearContents()
$ws.Range("AP4").Value = 2.172582619339046
$ws.Range("AQ4").Value = 18.34782608695652
